$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Batch No. cell (C3): replace the old fixed "BATCH44444" text with a
# date value, picking up the same date formatting/style already used by
# the "Required Date" cell (C4) right below it (clipboard format copy).
$ws.Range("C4").Copy()
$ws.Range("C3").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("C3").Value = 44537

# --- Item "Date in" column (D8:D17): switch from hard-coded dates to a
# live =TODAY() formula, same as already used for the "Required Date" cell.
$ws.Range("D8").Formula = "=TODAY()"
$ws.Range("D9").Formula = "=TODAY()"
$ws.Range("D10").Formula = "=TODAY()"
$ws.Range("D11").Formula = "=TODAY()"
$ws.Range("D12").Formula = "=TODAY()"
$ws.Range("D13").Formula = "=TODAY()"
$ws.Range("D14").Formula = "=TODAY()"
$ws.Range("D15").Formula = "=TODAY()"
$ws.Range("D16").Formula = "=TODAY()"
$ws.Range("D17").Formula = "=TODAY()"

# --- Move the active selection/clipboard focus from C4 to E5.
$ws.Range("E5").Select()
